# Finish the history data table
# - fill in the remaining two columns (Цел / Обхват) for the existing
#   "Балканска война" row (row 27)
# - append a new row (28) for "Първа световна"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- complete row 27 ("Балканска война") -----------------------------------
$ws.Range("E27").Value = "Териториално разширение и връщане на територия(не ги получава)"
$ws.Range("E27").WrapText = $true
$ws.Range("E27").VerticalAlignment = -4160   # xlTop

$ws.Range("F27").Value = "-"
$ws.Range("F27").VerticalAlignment = -4160   # xlTop

# --- add new row 28 ("Първа световна") --------------------------------------
$ws.Range("A28").Value = "Първа световна"
$ws.Range("A28").VerticalAlignment = -4160   # xlTop

$ws.Range("B28").Value = "1914-1918"
$ws.Range("B28").WrapText = $true
$ws.Range("B28").VerticalAlignment = -4160   # xlTop

$ws.Range("C28").Value = "Не"
$ws.Range("C28").VerticalAlignment = -4160   # xlTop

$ws.Range("D28").Value = "Централни сили(Германия, Италия, Австро-Унгария) vs Антантата"
$ws.Range("D28").WrapText = $true
$ws.Range("D28").VerticalAlignment = -4160   # xlTop

$ws.Range("E28").Value = "Териториално разширение и връщане на територия"
$ws.Range("E28").WrapText = $true
$ws.Range("E28").VerticalAlignment = -4160   # xlTop

$ws.Range("F28").Value = "-"
$ws.Range("F28").VerticalAlignment = -4160   # xlTop

# row 28 wraps to two lines, same as the similarly-shaped rows above it
$ws.Rows.Item(28).RowHeight = 28.8

# --- update view / selection to match the author's final state -------------
$ws.Range("B28").Select()
